$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column ("Schema reference"), shifting Part/Ref/Ammount/
# Price/Weight/Electrical-consumption one column to the right (A->B, B->C, ...).
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column.
$ws.Range("A1").Value = "Schema reference"

# Column widths: new column A gets a fresh width; column B (old column A)
# keeps its existing best-fit width untouched.
$ws.Columns.Item(1).ColumnWidth = 15.166666666666666

# Move/collapse the selection to A2 (previously E14).
$ws.Range("A2").Select() | Out-Null
